$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 4).Value = "43.284.36"
$ws.Cells.Item(2, 5).Value = "  +2.95%  "
$ws.Cells.Item(3, 4).Value = "2.314.85"
$ws.Cells.Item(3, 5).Value = "  +2.53%  "
$ws.Cells.Item(4, 4).Value = "'1.00"
$ws.Cells.Item(4, 4).Style = "Normal"
$ws.Cells.Item(4, 5).Value = "  -0.03%  "
$ws.Cells.Item(5, 4).Value = "'310.50"
$ws.Cells.Item(5, 4).Style = "Normal"
$ws.Cells.Item(5, 5).Value = "  +1.85%  "
$ws.Cells.Item(6, 4).Value = "'101.48"
$ws.Cells.Item(6, 4).Style = "Normal"
$ws.Cells.Item(6, 5).Value = "  +5.73%  "
$ws.Cells.Item(7, 4).Value = "'0.535"
$ws.Cells.Item(7, 4).Style = "Normal"
$ws.Cells.Item(7, 5).Value = "  +2.15%  "
$ws.Cells.Item(8, 5).Value = "  -0.04%  "
$ws.Cells.Item(9, 4).Value = "'0.525"
$ws.Cells.Item(9, 4).Style = "Normal"
$ws.Cells.Item(9, 5).Value = "  +7.27%  "
$ws.Cells.Item(10, 5).Value = "  +3.51%  "
$ws.Cells.Item(11, 5).Value = "  +3.58%  "
$ws.Cells.Item(12, 5).Value = "  +0.73%  "
$ws.Cells.Item(13, 4).Value = "'7.09"
$ws.Cells.Item(13, 4).Style = "Normal"
$ws.Cells.Item(13, 5).Value = "  +4.12%  "
$ws.Cells.Item(14, 4).Value = "2.673.76"
$ws.Cells.Item(14, 5).Value = "  +2.52%  "
$ws.Cells.Item(15, 4).Value = "'15.01"
$ws.Cells.Item(15, 4).Style = "Normal"
$ws.Cells.Item(15, 5).Value = "  +3.73%  "
$ws.Cells.Item(16, 4).Value = "2.312.76"
$ws.Cells.Item(16, 5).Value = "  +2.37%  "
$ws.Cells.Item(17, 4).Value = "'0.813"
$ws.Cells.Item(17, 4).Style = "Normal"
$ws.Cells.Item(17, 5).Value = "  +2.79%  "
$ws.Cells.Item(18, 4).Value = "43.180.77"
$ws.Cells.Item(18, 5).Value = "  +2.99%  "
$ws.Cells.Item(19, 4).Value = "'12.58"
$ws.Cells.Item(19, 4).Style = "Normal"
$ws.Cells.Item(19, 5).Value = "  +1.95%  "
$ws.Cells.Item(20, 5).Value = "  +2.25%  "
$ws.Cells.Item(21, 4).Value = "'6.15"
$ws.Cells.Item(21, 4).Style = "Normal"
$ws.Cells.Item(21, 5).Value = "  +3.33%  "
$ws.Cells.Item(22, 4).Value = "'68.45"
$ws.Cells.Item(22, 4).Style = "Normal"
$ws.Cells.Item(22, 5).Value = "  +0.04%  "
$ws.Cells.Item(23, 4).Value = "'241.39"
$ws.Cells.Item(23, 4).Style = "Normal"
$ws.Cells.Item(23, 5).Value = "  +1.84%  "
$ws.Cells.Item(24, 5).Value = "  +6.24%  "
$ws.Cells.Item(25, 5).Value = "  +3.52%  "
$ws.Cells.Item(26, 4).Value = "'1.00"
$ws.Cells.Item(26, 4).Style = "Normal"
$ws.Cells.Item(26, 5).Value = "  -0.02%  "
$ws.Cells.Item(27, 4).Value = "'24.71"
$ws.Cells.Item(27, 4).Style = "Normal"
$ws.Cells.Item(27, 5).Value = "  +5.10%  "
$ws.Cells.Item(28, 4).Value = "'37.55"
$ws.Cells.Item(28, 4).Style = "Normal"
$ws.Cells.Item(28, 5).Value = "  +3.12%  "
$ws.Cells.Item(29, 5).Value = "  +2.58%  "
$ws.Cells.Item(30, 5).Value = "  -0.16%  "
$ws.Cells.Item(31, 4).Value = "'167.76"
$ws.Cells.Item(31, 4).Style = "Normal"
$ws.Cells.Item(31, 5).Value = "  +4.31%  "
$ws.Cells.Item(32, 5).Value = "  +2.78%  "
$ws.Cells.Item(33, 4).Value = "'1.00"
$ws.Cells.Item(33, 4).Style = "Normal"
$ws.Cells.Item(33, 5).Value = "  -0.01%  "
$ws.Cells.Item(34, 4).Value = "'3.15"
$ws.Cells.Item(34, 4).Style = "Normal"
$ws.Cells.Item(34, 5).Value = "  -0.17%  "
$ws.Cells.Item(35, 5).Value = "  +5.74%  "
$ws.Cells.Item(36, 5).Value = "  +1.52%  "
$ws.Cells.Item(37, 5).Value = "  +3.01%  "
$ws.Cells.Item(38, 5).Value = "  +0.94%  "
$ws.Cells.Item(39, 4).Value = "'1.85"
$ws.Cells.Item(39, 4).Style = "Normal"
$ws.Cells.Item(39, 5).Value = "  +3.01%  "
$ws.Cells.Item(40, 5).Value = "  +2.26%  "
$ws.Cells.Item(41, 5).Value = "  +7.95%  "
$ws.Cells.Item(42, 4).Value = "'19.87"
$ws.Cells.Item(42, 4).Style = "Normal"
$ws.Cells.Item(42, 5).Value = "  +7.59%  "
$ws.Cells.Item(43, 4).Value = "'2.31"
$ws.Cells.Item(43, 4).Style = "Normal"
$ws.Cells.Item(43, 5).Value = "  +1.50%  "
$ws.Cells.Item(44, 5).Value = "  +3.69%  "
$ws.Cells.Item(45, 4).Value = "1.974.66"
$ws.Cells.Item(45, 5).Value = "  +0.77%  "
$ws.Cells.Item(46, 5).Value = "  +4.21%  "
$ws.Cells.Item(47, 4).Value = "'9.82"
$ws.Cells.Item(47, 4).Style = "Normal"
$ws.Cells.Item(47, 5).Value = "  -0.92%  "
$ws.Cells.Item(48, 5).Value = "  +19.19%  "
$ws.Cells.Item(49, 4).Value = "'55.82"
$ws.Cells.Item(49, 4).Style = "Normal"
$ws.Cells.Item(49, 5).Value = "  +4.85%  "
$ws.Cells.Item(50, 4).Value = "2.540.13"
$ws.Cells.Item(50, 5).Value = "  +2.41%  "
$ws.Cells.Item(51, 5).Value = "  +4.52%  "
